$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B,C to C,D)
$ws.Columns("B").Insert()

# Set header and value for the newly inserted column B
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_designation IN [''COTC007B'',''NCATS-COP01'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

# Match the wrap-text style used by A2 on the new B2 cell
$ws.Range("B2").WrapText = $true

# Match column B's width to column A's width
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Update the active selection to B2, matching the post-edit state
$ws.Range("B2").Select()
